# correction in sa algorithm and 746 logs
# Updates column C (Fitness) values for rows 2-252 on the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: start row, end row, new value (contiguous blocks of identical values)
$blocks = @(
    @(2,   3,   8856),
    @(4,   23,  8311),
    @(24,  63,  8159),
    @(64,  66,  7798),
    @(67,  88,  7343),
    @(172, 252, 7293)
)

foreach ($b in $blocks) {
    $startRow = $b[0]
    $endRow   = $b[1]
    $value    = $b[2]

    $rangeAddr = "C" + $startRow + ":C" + $endRow
    $ws.Range($rangeAddr).Value = $value
}
